# Fruta / hortaliza, semanal
# Insert a new weekly record as row 26 in "Fruta, Macroferia Regional de
# Talca - Chirimoya" (pushes the existing rows 26-31 down to 27-32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 26, shifting rows 26-31 -> 27-32
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly record
$ws.Range("A26").Value = 5
$ws.Range("B26").Value = "Macroferia Regional de Talca"
$ws.Range("C26").Value = "Maule"
$ws.Range("D26").Value = 44476
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100107
$ws.Range("H26").Value = "Otros"
$ws.Range("I26").Value = 100107002
$ws.Range("J26").Value = "Chirimoya"
$ws.Range("K26").Value = "Cultivar IV Región"
$ws.Range("L26").Value = "Especial"
$ws.Range("M26").Value = 100
$ws.Range("N26").Value = 30000
$ws.Range("O26").Value = 30000
$ws.Range("P26").Value = 30000
$ws.Range("Q26").Value = "$/bandeja 10 kilos"
$ws.Range("R26").Value = "Provincia de Limarí"
$ws.Range("S26").Value = 3000
$ws.Range("T26").Value = 10
